# Generate Report for Archive
#
# The localization-status report is regenerated: the in-progress status
# label changes from "Ready for handoff" to "In Translation" on every
# sheet that surfaces it, and the Status column on each sheet is resized
# to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status shown per-locale in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: status in column C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: status in column C
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
